$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Matriz_Resultados - a few win/loss/tie indicators were corrected
# from the previous (buggy) Diebold-Mariano comparison to the corrected one.
# ---------------------------------------------------------------------------
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("E2").Value = 0
$wsMatriz.Range("H2").Value = 0
$wsMatriz.Range("I2").Value = 0
$wsMatriz.Range("E3").Value = 0
$wsMatriz.Range("G4").Value = 0
$wsMatriz.Range("I4").Value = 0
$wsMatriz.Range("B5").Value = 0
$wsMatriz.Range("C5").Value = 0
$wsMatriz.Range("D7").Value = 0
$wsMatriz.Range("B8").Value = 0
$wsMatriz.Range("B9").Value = 0
$wsMatriz.Range("D9").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: P_valores - recomputed p-values (symmetric matrix, diagonal = 1)
# ---------------------------------------------------------------------------
$wsPval = $wb.Worksheets.Item("P_valores")

$wsPval.Range("C2").Value = 0.02137264465530153
$wsPval.Range("D2").Value = 0.03596467429554306
$wsPval.Range("E2").Value = 0.00423215505102692
$wsPval.Range("F2").Value = 0.00001551726171777901
$wsPval.Range("G2").Value = 0.00003411557354104922
$wsPval.Range("H2").Value = 0.005560685366990281
$wsPval.Range("I2").Value = 0.004033296139978493
$wsPval.Range("J2").Value = 0.00000000002491384876179836

$wsPval.Range("B3").Value = 0.02137264465530153
$wsPval.Range("D3").Value = 0.003664610293440118
$wsPval.Range("E3").Value = 0.002696502129224454
$wsPval.Range("F3").Value = 0.0001523575539801314
$wsPval.Range("G3").Value = 0.000007287232903552265
$wsPval.Range("H3").Value = 0.0005577210521019449
$wsPval.Range("I3").Value = 0.1719411947263698
$wsPval.Range("J3").Value = 0.000000000007242872968049596

$wsPval.Range("B4").Value = 0.03596467429554306
$wsPval.Range("C4").Value = 0.003664610293440118
$wsPval.Range("E4").Value = 0.01112372012711549
$wsPval.Range("F4").Value = 0.01272238286799432
$wsPval.Range("G4").Value = 0.001760689025062545
$wsPval.Range("H4").Value = 0.0855900837232495
$wsPval.Range("I4").Value = 0.003610530804192935
$wsPval.Range("J4").Value = 0.00000000001882050071344565

$wsPval.Range("B5").Value = 0.00423215505102692
$wsPval.Range("C5").Value = 0.002696502129224454
$wsPval.Range("D5").Value = 0.01112372012711549
$wsPval.Range("F5").Value = 0.6496163373296979
$wsPval.Range("G5").Value = 0.4690370762678189
$wsPval.Range("H5").Value = 0.8239416967008026
$wsPval.Range("I5").Value = 0.0009768522151367698
$wsPval.Range("J5").Value = 0.00000437413006637577

$wsPval.Range("B6").Value = 0.00001551726171777901
$wsPval.Range("C6").Value = 0.0001523575539801314
$wsPval.Range("D6").Value = 0.01272238286799432
$wsPval.Range("E6").Value = 0.6496163373296979
$wsPval.Range("G6").Value = 0.08980689798735786
$wsPval.Range("H6").Value = 0.9184222965522681
$wsPval.Range("I6").Value = 0.00004938113564323388
$wsPval.Range("J6").Value = 0.0000001109920231279204

$wsPval.Range("B7").Value = 0.00003411557354104922
$wsPval.Range("C7").Value = 0.000007287232903552265
$wsPval.Range("D7").Value = 0.001760689025062545
$wsPval.Range("E7").Value = 0.4690370762678189
$wsPval.Range("F7").Value = 0.08980689798735786
$wsPval.Range("H7").Value = 0.02053841934637957
$wsPval.Range("I7").Value = 0.00001118578073278087
$wsPval.Range("J7").Value = 0.0000008649737659460754

$wsPval.Range("B8").Value = 0.005560685366990281
$wsPval.Range("C8").Value = 0.0005577210521019449
$wsPval.Range("D8").Value = 0.0855900837232495
$wsPval.Range("E8").Value = 0.8239416967008026
$wsPval.Range("F8").Value = 0.9184222965522681
$wsPval.Range("G8").Value = 0.02053841934637957
$wsPval.Range("I8").Value = 0.0002021295352028218
$wsPval.Range("J8").Value = 0.000001852158119541869

$wsPval.Range("B9").Value = 0.004033296139978493
$wsPval.Range("C9").Value = 0.1719411947263698
$wsPval.Range("D9").Value = 0.003610530804192935
$wsPval.Range("E9").Value = 0.0009768522151367698
$wsPval.Range("F9").Value = 0.00004938113564323388
$wsPval.Range("G9").Value = 0.00001118578073278087
$wsPval.Range("H9").Value = 0.0002021295352028218
$wsPval.Range("J9").Value = 0.0000000007305507310206849

$wsPval.Range("B10").Value = 0.00000000002491384876179836
$wsPval.Range("C10").Value = 0.000000000007242872968049596
$wsPval.Range("D10").Value = 0.00000000001882050071344565
$wsPval.Range("E10").Value = 0.00000437413006637577
$wsPval.Range("F10").Value = 0.0000001109920231279204
$wsPval.Range("G10").Value = 0.0000008649737659460754
$wsPval.Range("H10").Value = 0.000001852158119541869
$wsPval.Range("I10").Value = 0.0000000007305507310206849

# ---------------------------------------------------------------------------
# Sheet 3: Estadisticos_DM - recomputed Diebold-Mariano statistics (antisymmetric)
# ---------------------------------------------------------------------------
$wsStat = $wb.Worksheets.Item("Estadisticos_DM")

$wsStat.Range("C2").Value = 2.590404370765578
$wsStat.Range("D2").Value = -2.31991588709203
$wsStat.Range("E2").Value = -3.409500173225791
$wsStat.Range("F2").Value = -6.438593383659957
$wsStat.Range("G2").Value = -5.972667226209714
$wsStat.Range("H2").Value = -3.272278828574468
$wsStat.Range("I2").Value = 3.433698227646866
$wsStat.Range("J2").Value = -18.79564760716249

$wsStat.Range("B3").Value = -2.590404370765578
$wsStat.Range("D3").Value = -3.481912368932113
$wsStat.Range("E3").Value = -3.63641667631377
$wsStat.Range("F3").Value = -5.132333516325322
$wsStat.Range("G3").Value = -6.90316674856979
$wsStat.Range("H3").Value = -4.442167822176491
$wsStat.Range("I3").Value = 1.439685335625883
$wsStat.Range("J3").Value = -20.59111536364859

$wsStat.Range("B4").Value = 2.31991588709203
$wsStat.Range("C4").Value = 3.481912368932113
$wsStat.Range("E4").Value = -2.923038319617348
$wsStat.Range("F4").Value = -2.855056884245833
$wsStat.Range("G4").Value = -3.851971882875567
$wsStat.Range("H4").Value = -1.849615575723482
$wsStat.Range("I4").Value = 3.4893921035369
$wsStat.Range("J4").Value = -19.18996178529352

$wsStat.Range("B5").Value = 3.409500173225791
$wsStat.Range("C5").Value = 3.63641667631377
$wsStat.Range("D5").Value = 2.923038319617348
$wsStat.Range("F5").Value = 0.4642345206799753
$wsStat.Range("G5").Value = -0.7442448084827762
$wsStat.Range("H5").Value = 0.2266895502107374
$wsStat.Range("I5").Value = 4.152470992738388
$wsStat.Range("J5").Value = -7.227460767262928

$wsStat.Range("B6").Value = 6.438593383659957
$wsStat.Range("C6").Value = 5.132333516325322
$wsStat.Range("D6").Value = 2.855056884245833
$wsStat.Range("E6").Value = -0.4642345206799753
$wsStat.Range("G6").Value = -1.822482023640444
$wsStat.Range("H6").Value = -0.1042852008596833
$wsStat.Range("I6").Value = 5.759867461460644
$wsStat.Range("J6").Value = -9.861546994881959

$wsStat.Range("B7").Value = 5.972667226209714
$wsStat.Range("C7").Value = 6.90316674856979
$wsStat.Range("D7").Value = 3.851971882875567
$wsStat.Range("E7").Value = 0.7442448084827762
$wsStat.Range("F7").Value = 1.822482023640444
$wsStat.Range("H7").Value = 2.61086072006762
$wsStat.Range("I7").Value = 6.637565897076734
$wsStat.Range("J7").Value = -8.319983111453229

$wsStat.Range("B8").Value = 3.272278828574468
$wsStat.Range("C8").Value = 4.442167822176491
$wsStat.Range("D8").Value = 1.849615575723482
$wsStat.Range("E8").Value = -0.2266895502107374
$wsStat.Range("F8").Value = 0.1042852008596833
$wsStat.Range("G8").Value = -2.61086072006762
$wsStat.Range("I8").Value = 4.979308767983962
$wsStat.Range("J8").Value = -7.794274463586601

$wsStat.Range("B9").Value = -3.433698227646866
$wsStat.Range("C9").Value = -1.439685335625883
$wsStat.Range("D9").Value = -3.4893921035369
$wsStat.Range("E9").Value = -4.152470992738388
$wsStat.Range("F9").Value = -5.759867461460644
$wsStat.Range("G9").Value = -6.637565897076734
$wsStat.Range("H9").Value = -4.979308767983962
$wsStat.Range("J9").Value = -14.59514854873765

$wsStat.Range("B10").Value = 18.79564760716249
$wsStat.Range("C10").Value = 20.59111536364859
$wsStat.Range("D10").Value = 19.18996178529352
$wsStat.Range("E10").Value = 7.227460767262928
$wsStat.Range("F10").Value = 9.861546994881959
$wsStat.Range("G10").Value = 8.319983111453229
$wsStat.Range("H10").Value = 7.794274463586601
$wsStat.Range("I10").Value = 14.59514854873765

# ---------------------------------------------------------------------------
# Sheet 4: Resumen - summary table re-sorted/recomputed from the corrected
# win/loss/tie matrix (model names stay matched to their own ECRPS_Medio).
# ---------------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")

$wsResumen.Range("A2").Value = "DeepAR"
$wsResumen.Range("B2").Value = 5
$wsResumen.Range("C2").Value = 0
$wsResumen.Range("D2").Value = 3
$wsResumen.Range("E2").Value = 62.5
$wsResumen.Range("F2").Value = 0.6266411458939485

$wsResumen.Range("A3").Value = "Sieve Bootstrap"
$wsResumen.Range("B3").Value = 4
$wsResumen.Range("C3").Value = 0
$wsResumen.Range("D3").Value = 4
$wsResumen.Range("E3").Value = 50
$wsResumen.Range("F3").Value = 0.6404604772955071

$wsResumen.Range("A4").Value = "Block Bootstrapping"
$wsResumen.Range("B4").Value = 3
$wsResumen.Range("C4").Value = 0
$wsResumen.Range("D4").Value = 5
$wsResumen.Range("E4").Value = 37.5
$wsResumen.Range("F4").Value = 0.6615292813969161

$wsResumen.Range("A5").Value = "LSPMW"
$wsResumen.Range("B5").Value = 1
$wsResumen.Range("C5").Value = 1
$wsResumen.Range("D5").Value = 6
$wsResumen.Range("E5").Value = 12.5
$wsResumen.Range("F5").Value = 0.7240299305229702

$wsResumen.Range("A6").Value = "LSPM"
$wsResumen.Range("B6").Value = 1
$wsResumen.Range("C6").Value = 0
$wsResumen.Range("D6").Value = 7
$wsResumen.Range("E6").Value = 12.5
$wsResumen.Range("F6").Value = 0.6791407632310854

$wsResumen.Range("A7").Value = "AREPD"
$wsResumen.Range("B7").Value = 1
$wsResumen.Range("C7").Value = 3
$wsResumen.Range("D7").Value = 4
$wsResumen.Range("E7").Value = 12.5
$wsResumen.Range("F7").Value = 0.7156707932082101

$wsResumen.Range("A8").Value = "MCPS"
$wsResumen.Range("B8").Value = 1
$wsResumen.Range("C8").Value = 3
$wsResumen.Range("D8").Value = 4
$wsResumen.Range("E8").Value = 12.5
$wsResumen.Range("F8").Value = 0.7447845475734969

$wsResumen.Range("A9").Value = "AV-MCPS"
$wsResumen.Range("B9").Value = 1
$wsResumen.Range("C9").Value = 2
$wsResumen.Range("D9").Value = 5
$wsResumen.Range("E9").Value = 12.5
$wsResumen.Range("F9").Value = 0.7175246340748288

$wsResumen.Range("A10").Value = "EnCQR-LSTM"
$wsResumen.Range("B10").Value = 0
$wsResumen.Range("C10").Value = 8
$wsResumen.Range("D10").Value = 0
$wsResumen.Range("E10").Value = 0
$wsResumen.Range("F10").Value = 0.8972757458314203
